$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

# Append the new "map_source" translation row (key / French / English) that
# was missing from the sheet -- fixes the issue reported via mail.
$ws.Range("A54").Value = "map_source"
$ws.Range("C54").Value = "source: Bibliothèque de l'Institut national d'histoire de l'art, collections Jacques Doucet,* 12 RES 870*, crédit photo : INHA"
$ws.Range("D54").Value = "image source: Bibliothèque de l'Institut national d'histoire de l'art, collections Jacques Doucet, 12 RES 870, photo credit:INHA"

# This is a long wrapped string, same as the other "source" style rows --
# give it the same taller row height those use.
$ws.Rows.Item(54).RowHeight = 46.25

# Scroll the frozen view down a bit and leave the new row's first cell
# selected, same as when the row was added by hand.
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$null = $ws.Range("A54").Select()
